$wb = $excel.ActiveWorkbook
$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# Overview sheet
$ov.Range("A1").Value = "File Name"
$ov.Range("B1").Value = "zh-cn"
$ov.Range("C1").Value = "de-de"
$ov.Range("A2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("A3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

# zh-cn sheet
$zh.Range("A1").Value = "Source File Name"
$zh.Range("B1").Value = "Status"
$zh.Range("C1").Value = "Latest Handoff File"
$zh.Range("D1").Value = "Latest Handoff Datetime"
$zh.Range("E1").Value = "Latest Target File"
$zh.Range("F1").Value = "Latest Handback File"
$zh.Range("G1").Value = "Latest Handback DateTime"
$zh.Range("H1").Value = "Handoff Reason"
$zh.Range("I1").Value = "Dependency From"
$zh.Range("A2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
$zh.Range("D2").Value = "2016-01-17 03:18:19"
$zh.Range("E2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$zh.Range("F2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.zh-cn.xlf"
$zh.Range("G2").Value = "2016-01-17 03:19:00"
$zh.Range("H2").Value = "Include"
$zh.Range("A3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
$zh.Range("D3").Value = "2016-01-17 03:20:01"
$zh.Range("E3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$zh.Range("F3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.zh-cn.xlf"
$zh.Range("G3").Value = "2016-01-17 03:19:00"
$zh.Range("H3").Value = "Include"
$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

# de-de sheet
$de.Range("A1").Value = "Source File Name"
$de.Range("B1").Value = "Status"
$de.Range("C1").Value = "Latest Handoff File"
$de.Range("D1").Value = "Latest Handoff Datetime"
$de.Range("E1").Value = "Latest Target File"
$de.Range("F1").Value = "Latest Handback File"
$de.Range("G1").Value = "Latest Handback DateTime"
$de.Range("H1").Value = "Handoff Reason"
$de.Range("I1").Value = "Dependency From"
$de.Range("A2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
$de.Range("D2").Value = "2016-01-17 03:18:30"
$de.Range("E2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.md"
$de.Range("F2").Value = "a1eefbe0-7ce3-406d-afc9-ca13d215af41.b2d1c755cfdc5700f172a16459738075dc1fff73.de-de.xlf"
$de.Range("G2").Value = "2016-01-17 03:19:17"
$de.Range("H2").Value = "Include"
$de.Range("A3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
$de.Range("D3").Value = "2016-01-17 03:20:11"
$de.Range("E3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.md"
$de.Range("F3").Value = "3120aef5-5742-44c4-bfc8-f48d3381e7be.787c48ea0e98c471fb60715da191c4fce12d627f.de-de.xlf"
$de.Range("G3").Value = "2016-01-17 03:19:17"
$de.Range("H3").Value = "Include"
$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"
